$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 190: LeetCode 2048 - Next Greater Numerically Balanced Number
$ws.Range("A190").Value = 2048
$ws.Range("B190").Value = "Next Greater Numerically Balanced Number"
$ws.Range("C190").Value = "#math #emueration "
$ws.Range("D190").Value = "medium"
$ws.Range("E190").Value = 0
$ws.Range("F190").Value = 1
$ws.Range("G190").Value = 25
$ws.Range("H190").Value = 45954
$ws.Range("H190").NumberFormat = "m/d/yy"
$ws.Range("I190").Value = 45954
$ws.Range("I190").NumberFormat = "m/d/yy"
$ws.Rows.Item(190).RowHeight = 51

# Row 191: LeetCode 1716 - Calculate Money in Leetcode Bank
$ws.Range("A191").Value = 1716
$ws.Range("B191").Value = "Calculate Money in Leetcode Bank"
$ws.Range("C191").Value = "#math"
$ws.Range("D191").Value = "easy"
$ws.Range("E191").Value = 1
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 9
$ws.Range("H191").Value = 45955
$ws.Range("H191").NumberFormat = "m/d/yy"
$ws.Range("I191").Value = 45955
$ws.Range("I191").NumberFormat = "m/d/yy"
$ws.Rows.Item(191).RowHeight = 34

# Row 192: LeetCode 2125 - Number of Laser Beams in a Bank
$ws.Range("A192").Value = 2125
$ws.Range("B192").Value = "Number of Laser Beams in a Bank"
$ws.Range("C192").Value = "#string  #greedy "
$ws.Range("D192").Value = "medium"
$ws.Range("E192").Value = 1
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 8
$ws.Range("H192").Value = 45957
$ws.Range("H192").NumberFormat = "m/d/yy"
$ws.Range("I192").Value = 45957
$ws.Range("I192").NumberFormat = "m/d/yy"
$ws.Rows.Item(192).RowHeight = 34

# Update view state to match where the author ended up after the edits
$ws.Range("E191").Select() | Out-Null
